$d = $word.ActiveDocument

# --- Step 1: Replace the whole span of old "demonstration description" content
#     (original paragraphs 4-8, 1-indexed) with the new merged text. We target
#     a range that begins at the start of paragraph 4 and ends just before the
#     paragraph mark that terminates paragraph 8 (leaving that paragraph mark,
#     and the following paragraphs, alone) so the assignment lands cleanly in
#     a single run instead of tripping the "assign across a pilcrow" quirk.

$p4 = $d.Paragraphs.Item(4)
$p8 = $d.Paragraphs.Item(8)

$newText = "A system comprised of a laptop, cellular modem, and a load-box will be used to demonstrate the functionality of the battery power monitor (BPM). The BPM will monitor the voltage and current of the cellular modem during transmission and log all data to an on-board SD card. In addition to data logging, minimum, maximum, and average values of the power consumption will be displayed on an attached LCD display. These statistical calculations can be reset via a push-button.  Real-time sensor data will be sent over USB to the laptop and displayed on screen. In addition to monitoring the cellular modem, the general nature of the BPM will be demonstrated by using a load box to simulate a variety of different loads.  "

$full = $d.Range($p4.Range.Start, $p8.Range.End - 1)
$full.Text = $newText

# --- Step 2: the old paragraphs 5-8 are now empty orphans sitting right after
#     the (now rewritten) paragraph 4; delete them so paragraph 4 again directly
#     precedes the trailing " " paragraph.

for ($k = 0; $k -lt 4; $k++) {
    $d.Paragraphs.Item(5).Range.Delete()
}

# --- Step 3: restore the "_GoBack" bookmark at its new home - right after
#     "...on-board SD card. " and before "In addition to data logging...".

$p4 = $d.Paragraphs.Item(4)
$marker = "an on-board SD card. "
$bmPos = $p4.Range.Start + $newText.IndexOf($marker) + $marker.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
